# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.981.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +6.64%  '

# Row 3
$ws.Range("D3").Value = "'1.740.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.09%  '

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").Value = "'228.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.37%  '

# Row 6
$ws.Range("D6").Value = "'0.5464"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.26%  '

# Row 7
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.18%  '

# Row 8
$ws.Range("D8").Value = "'0.2779"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.29%  '

# Row 9
$ws.Range("D9").Value = "'0.06728"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.79%  '

# Row 10
$ws.Range("D10").Value = "'21.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.12%  '

# Row 11
$ws.Range("D11").Value = "'0.07796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.89%  '

# Row 12
$ws.Range("D12").Value = "'4.708"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.38%  '

# Row 13
$ws.Range("D13").Value = "'1.979.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.06%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = "'1.714.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.28%  '

# Row 15
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = "'0.6018"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.83%  '

# Row 16
$ws.Range("D16").Value = "'0.0₅8439"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.23%  '

# Row 17
$ws.Range("D17").Value = "'69.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.19%  '

# Row 18
$ws.Range("D18").Value = "'27.973.19"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'226.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +17.55%  '

# Row 20
$ws.Range("D20").Value = "'4.849"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.21%  '

# Row 22
$ws.Range("D22").Value = "'10.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.18%  '

# Row 23
$ws.Range("D23").Value = "'6.242"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.93%  '

# Row 24
$ws.Range("D24").Value = "'1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '

# Row 25
$ws.Range("D25").Value = "'146.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.69%  '

# Row 26
$ws.Range("E26").Value = '  +4.25%  '

# Row 27
$ws.Range("D27").Value = "'7.478"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.79%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = "'1.663"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.86%  '

# Row 29
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = "'17.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.34%  '

# Row 30
$ws.Range("D30").Value = "'0.05698"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.06%  '

# Row 31
$ws.Range("D31").Value = "'1.316"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.92%  '

# Row 32
$ws.Range("D32").Value = "'3.705"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.50%  '

# Row 33
$ws.Range("D33").Value = "'3.541"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.39%  '

# Row 34
$ws.Range("D34").Value = "'1.670"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.47%  '

# Row 35
$ws.Range("D35").Value = "'0.9851"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.16%  '

# Row 36
$ws.Range("D36").Value = "'2.863"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.13%  '

# Row 37
$ws.Range("D37").Value = "'2.454"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.69%  '

# Row 38
$ws.Range("D38").Value = "'0.5961"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.49%  '

# Row 39
$ws.Range("D39").Value = "'0.01675"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.65%  '

# Row 40
$ws.Range("E40").Value = '  +0.25%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'0.8497"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.73%  '

# Row 42
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = "'1.048.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.78%  '

# Row 43
$ws.Range("E43").Value = '  -0.13%  '

# Row 44
$ws.Range("D44").Value = "'102.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.31%  '

# Row 45
$ws.Range("D45").Value = "'1.884.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.99%  '

# Row 46
$ws.Range("E46").Value = '  +13.00%  '

# Row 47
$ws.Range("D47").Value = "'60.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.12%  '

# Row 48
$ws.Range("D48").Value = "'8.345"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.05%  '

# Row 49
$ws.Range("D49").Value = "'1.011"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.40%  '

# Row 50
$ws.Range("D50").Value = "'0.4427"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.80%  '

# Row 51
$ws.Range("D51").Value = "'0.05320"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.56%  '
